$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C402").Value = 2
